# post_test.xlsx — add three new per-row columns (L/M/O) that capture the
# size overrides for the left/right "up" images plus the allowed response
# keys for each trial, shifting the old "show_slider" column from L to N.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the existing L column ("show_slider") into the new N column ---
$ws.Range("N1").Value = "show_slider"
$ws.Range("N2").Value = " "
$ws.Range("N3").Value = " "
$ws.Range("N4").Value = " "
$ws.Range("N5").Value = " "
$ws.Range("N6").Value = " "
$ws.Range("N7").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("N9").Value = 0
$ws.Range("N10").Value = 0

# --- New column O: allowed_keys ---
$ws.Range("O1").Value = "allowed_keys"
$ws.Range("O2").Value = "('x', 'm')"
$ws.Range("O3").Value = "('x', 'm')"
$ws.Range("O4").Value = "('x', 'm')"
$ws.Range("O5").Value = "('x', 'm')"
$ws.Range("O6").Value = "('1', '2','3', '4')"
$ws.Range("O7").Value = "('space')"
$ws.Range("O8").Value = "('space')"
$ws.Range("O9").Value = "('space')"
$ws.Range("O10").Value = "('space')"

# --- New column headers L/M: size_left_up_img / size_right_img_up ---
$ws.Range("L1").Value = "size_left_up_img"
$ws.Range("M1").Value = "size_right_img_up"

# --- Existing K column: two rows get the new explicit size ---
$ws.Range("K2").Value = "[0.35,0.25]"
$ws.Range("K3").Value = "[0.35,0.25]"

# --- New columns L/M data rows ---
$ws.Range("L2").Value = 0.3
$ws.Range("M2").Value = 0.3
$ws.Range("L3").Value = 0.3
$ws.Range("M3").Value = 0.3
$ws.Range("L4").Value = "[0.35,0.25]"
$ws.Range("M4").Value = 0.3
$ws.Range("L5").Value = "[0.35,0.25]"
$ws.Range("M5").Value = 0.3
$ws.Range("L6").Value = 0.2
$ws.Range("M6").Value = 0.2
$ws.Range("L7").Value = 0.3
$ws.Range("M7").Value = "[0.35,0.25]"
$ws.Range("L8").Value = 0.3
$ws.Range("M8").Value = "[0.35,0.25]"
$ws.Range("L9").Value = 0.3
$ws.Range("M9").Value = "[0.35,0.25]"
$ws.Range("L10").Value = 0.3
$ws.Range("M10").Value = "[0.35,0.25]"

# --- Final cursor position left by the editor ---
$ws.Range("Q5").Select()
